# Add files via upload
# Update the two transaction description cells: bump the year from 2024 to
# 2025 and normalize the date/description separator to " - ".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "13.02.2025 - Otistics Kel'el Ware karşılığında Out of Po'ya 2 Dolar vermiştir. (395-404)"
$ws.Range("A12").Value = "19.02.2025 - Otistics Portis karşılığında (ya da dize yatırma da diyebiliriz) Los Yahoo'ya 1 Dolar vermiştir. (394-400)"

# Move the active selection to A12 (was C5).
$ws.Range("A12").Select()
